$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 952 (pushes the existing 952:976 block down to 958:982,
# matching the new sheet dimension A1:R982). Insert() on a full-row range copies the
# formatting of the row above, which already carries the date style (s="2") on column D.
$ws.Rows("952:957").Insert()

# New weekly price rows (date 2021-09-09 = serial 44448) for "Zafiro rojo" / "Zafiro verde".
# Columns A,B,C,E,F,G,R are constant across this whole sub-block (same market/product).

# Row 952: Zafiro rojo / Extra
$ws.Cells.Item(952,1).Value = 8
$ws.Cells.Item(952,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(952,3).Value = "Coquimbo"
$ws.Cells.Item(952,4).Value = 44448
$ws.Cells.Item(952,5).Value = 4
$ws.Cells.Item(952,6).Value = 100112002
$ws.Cells.Item(952,7).Value = "Pimiento"
$ws.Cells.Item(952,8).Value = "Zafiro rojo"
$ws.Cells.Item(952,9).Value = "Extra"
$ws.Cells.Item(952,10).Value = 760
$ws.Cells.Item(952,11).Value = 39000
$ws.Cells.Item(952,12).Value = 40000
$ws.Cells.Item(952,13).Value = 39500
$ws.Cells.Item(952,14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(952,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(952,16).Value = 2633
$ws.Cells.Item(952,17).Value = 15
$ws.Cells.Item(952,18).Value = "Hortaliza"

# Row 953: Zafiro rojo / Primera
$ws.Cells.Item(953,1).Value = 8
$ws.Cells.Item(953,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(953,3).Value = "Coquimbo"
$ws.Cells.Item(953,4).Value = 44448
$ws.Cells.Item(953,5).Value = 4
$ws.Cells.Item(953,6).Value = 100112002
$ws.Cells.Item(953,7).Value = "Pimiento"
$ws.Cells.Item(953,8).Value = "Zafiro rojo"
$ws.Cells.Item(953,9).Value = "Primera"
$ws.Cells.Item(953,10).Value = 480
$ws.Cells.Item(953,11).Value = 37000
$ws.Cells.Item(953,12).Value = 38000
$ws.Cells.Item(953,13).Value = 37500
$ws.Cells.Item(953,14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(953,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(953,16).Value = 2500
$ws.Cells.Item(953,17).Value = 15
$ws.Cells.Item(953,18).Value = "Hortaliza"

# Row 954: Zafiro rojo / Segunda
$ws.Cells.Item(954,1).Value = 8
$ws.Cells.Item(954,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(954,3).Value = "Coquimbo"
$ws.Cells.Item(954,4).Value = 44448
$ws.Cells.Item(954,5).Value = 4
$ws.Cells.Item(954,6).Value = 100112002
$ws.Cells.Item(954,7).Value = "Pimiento"
$ws.Cells.Item(954,8).Value = "Zafiro rojo"
$ws.Cells.Item(954,9).Value = "Segunda"
$ws.Cells.Item(954,10).Value = 360
$ws.Cells.Item(954,11).Value = 35000
$ws.Cells.Item(954,12).Value = 36000
$ws.Cells.Item(954,13).Value = 35500
$ws.Cells.Item(954,14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(954,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(954,16).Value = 2367
$ws.Cells.Item(954,17).Value = 15
$ws.Cells.Item(954,18).Value = "Hortaliza"

# Row 955: Zafiro verde / Extra
$ws.Cells.Item(955,1).Value = 8
$ws.Cells.Item(955,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(955,3).Value = "Coquimbo"
$ws.Cells.Item(955,4).Value = 44448
$ws.Cells.Item(955,5).Value = 4
$ws.Cells.Item(955,6).Value = 100112002
$ws.Cells.Item(955,7).Value = "Pimiento"
$ws.Cells.Item(955,8).Value = "Zafiro verde"
$ws.Cells.Item(955,9).Value = "Extra"
$ws.Cells.Item(955,10).Value = 640
$ws.Cells.Item(955,11).Value = 34000
$ws.Cells.Item(955,12).Value = 35000
$ws.Cells.Item(955,13).Value = 34500
$ws.Cells.Item(955,14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(955,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(955,16).Value = 2300
$ws.Cells.Item(955,17).Value = 15
$ws.Cells.Item(955,18).Value = "Hortaliza"

# Row 956: Zafiro verde / Primera
$ws.Cells.Item(956,1).Value = 8
$ws.Cells.Item(956,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(956,3).Value = "Coquimbo"
$ws.Cells.Item(956,4).Value = 44448
$ws.Cells.Item(956,5).Value = 4
$ws.Cells.Item(956,6).Value = 100112002
$ws.Cells.Item(956,7).Value = "Pimiento"
$ws.Cells.Item(956,8).Value = "Zafiro verde"
$ws.Cells.Item(956,9).Value = "Primera"
$ws.Cells.Item(956,10).Value = 400
$ws.Cells.Item(956,11).Value = 32000
$ws.Cells.Item(956,12).Value = 33000
$ws.Cells.Item(956,13).Value = 32500
$ws.Cells.Item(956,14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(956,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(956,16).Value = 2167
$ws.Cells.Item(956,17).Value = 15
$ws.Cells.Item(956,18).Value = "Hortaliza"

# Row 957: Zafiro verde / Segunda
$ws.Cells.Item(957,1).Value = 8
$ws.Cells.Item(957,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(957,3).Value = "Coquimbo"
$ws.Cells.Item(957,4).Value = 44448
$ws.Cells.Item(957,5).Value = 4
$ws.Cells.Item(957,6).Value = 100112002
$ws.Cells.Item(957,7).Value = "Pimiento"
$ws.Cells.Item(957,8).Value = "Zafiro verde"
$ws.Cells.Item(957,9).Value = "Segunda"
$ws.Cells.Item(957,10).Value = 300
$ws.Cells.Item(957,11).Value = 28000
$ws.Cells.Item(957,12).Value = 29000
$ws.Cells.Item(957,13).Value = 28500
$ws.Cells.Item(957,14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(957,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(957,16).Value = 1900
$ws.Cells.Item(957,17).Value = 15
$ws.Cells.Item(957,18).Value = "Hortaliza"
